# TC29_VerifyALL_Links_Myaccount.xlsx - "Logic change for Logged in User"
#
# The login sequence used to be:
#   Row3: CLICK  LoginOption
#   Row4: ENTERTEXT Uname1 / Uname
#   Row5: ENTERTEXT Password1 / Password
#   Row6: CLICK  LoginButton1
#   Row7: WAIT
#   Row8: VERIFY_WEBELEMENT_PRESENT Welcomeelement
#
# It becomes:
#   Row3: CLICK  LoginOption
#   Row4: CLICK  LoginURL            <-- NEW ROW
#   Row5: ENTERTEXT Uname / Uname
#   Row6: ENTERTEXT Password / Password
#   Row7: CLICK  LoginButton
#   Row8: WAIT
#   Row9: VERIFY_WEBELEMENT_PRESENT Welcomeelement
#
# i.e. a new row is inserted right after row 3, and the old
# Uname1/Password1/LoginButton1 object names become the generic
# Uname/Password/LoginButton names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 4, shifting everything below down by one.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 (CLICK on LoginURL, CSS object type).
$ws.Range("B4").Value = "CLICK"
$ws.Range("C4").Value = "LoginURL"
$ws.Range("D4").Value = "CSS"

# Copy formatting from the row above (row 3) so the new row matches the
# look & feel of the rest of the table.
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Make sure the values are correct after the paste-special (paste format
# only should have kept formulas/values intact, but set again to be safe).
$ws.Range("B4").Value = "CLICK"
$ws.Range("C4").Value = "LoginURL"
$ws.Range("D4").Value = "CSS"
$ws.Range("E4").Value = ""

# Rename the old Uname1 / Password1 / LoginButton1 object identifiers
# (now living on rows 5, 6 and 7 after the insert) to the generic
# Uname / Password / LoginButton names used going forward.
$ws.Range("C5").Value = "Uname"
$ws.Range("C6").Value = "Password"
$ws.Range("C7").Value = "LoginButton"

# Restore the default view (Excel had scrolled/selected near the bottom of
# the sheet; reset it back to the top of the table).
$ws.Range("C5:C7").Select()

$wb.Save()
